$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2023-08-02 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-08-03 Thursday", 2) | Out-Null

# Update table cell values (direct cell addressing avoids text-collision issues)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "49×32=1568"
$t.Cell(1, 2).Range.Text = "97×16=1552"
$t.Cell(1, 3).Range.Text = "84×66=5544"
$t.Cell(1, 4).Range.Text = "100×53=5300"
$t.Cell(1, 5).Range.Text = "93×56=5208"
$t.Cell(2, 1).Range.Text = "93×92=8556"
$t.Cell(2, 2).Range.Text = "67×82=5494"
$t.Cell(2, 3).Range.Text = "12×38=456"
$t.Cell(2, 4).Range.Text = "56×34=1904"
$t.Cell(2, 5).Range.Text = "15×34=510"
$t.Cell(3, 1).Range.Text = "70×62=4340"
$t.Cell(3, 2).Range.Text = "98×35=3430"
$t.Cell(3, 3).Range.Text = "59×85=5015"
$t.Cell(3, 4).Range.Text = "35×29=1015"
$t.Cell(3, 5).Range.Text = "93×46=4278"
$t.Cell(4, 1).Range.Text = "57×11=627"
$t.Cell(4, 2).Range.Text = "14×68=952"
$t.Cell(4, 3).Range.Text = "17×84=1428"
$t.Cell(4, 4).Range.Text = "70×15=1050"
$t.Cell(4, 5).Range.Text = "86×73=6278"
$t.Cell(5, 1).Range.Text = "10×79=790"
$t.Cell(5, 2).Range.Text = "72×69=4968"
$t.Cell(5, 3).Range.Text = "94×62=5828"
$t.Cell(5, 4).Range.Text = "60×87=5220"
$t.Cell(5, 5).Range.Text = "76×33=2508"
$t.Cell(6, 1).Range.Text = "78×28=2184"
$t.Cell(6, 2).Range.Text = "88×89=7832"
$t.Cell(6, 3).Range.Text = "54×66=3564"
$t.Cell(6, 4).Range.Text = "17×66=1122"
$t.Cell(6, 5).Range.Text = "46×63=2898"
$t.Cell(7, 1).Range.Text = "87×59=5133"
$t.Cell(7, 2).Range.Text = "99×93=9207"
$t.Cell(7, 3).Range.Text = "84×81=6804"
$t.Cell(7, 4).Range.Text = "15×99=1485"
$t.Cell(7, 5).Range.Text = "13×34=442"
$t.Cell(8, 1).Range.Text = "88×99=8712"
$t.Cell(8, 2).Range.Text = "39×16=624"
$t.Cell(8, 3).Range.Text = "52×80=4160"
$t.Cell(8, 4).Range.Text = "86×36=3096"
$t.Cell(8, 5).Range.Text = "55×55=3025"
$t.Cell(9, 1).Range.Text = "94×69=6486"
$t.Cell(9, 2).Range.Text = "41×13=533"
$t.Cell(9, 3).Range.Text = "97×70=6790"
$t.Cell(9, 4).Range.Text = "41×63=2583"
$t.Cell(9, 5).Range.Text = "81×73=5913"
$t.Cell(10, 1).Range.Text = "78×95=7410"
$t.Cell(10, 2).Range.Text = "83×41=3403"
$t.Cell(10, 3).Range.Text = "34×36=1224"
$t.Cell(10, 4).Range.Text = "37×52=1924"
$t.Cell(10, 5).Range.Text = "78×29=2262"
$t.Cell(11, 1).Range.Text = "62×51=3162"
$t.Cell(11, 2).Range.Text = "46×60=2760"
$t.Cell(11, 3).Range.Text = "51×59=3009"
$t.Cell(11, 4).Range.Text = "20×80=1600"
$t.Cell(11, 5).Range.Text = "67×26=1742"
$t.Cell(12, 1).Range.Text = "66×61=4026"
$t.Cell(12, 2).Range.Text = "98×92=9016"
$t.Cell(12, 3).Range.Text = "12×16=192"
$t.Cell(12, 4).Range.Text = "73×81=5913"
$t.Cell(12, 5).Range.Text = "97×38=3686"
$t.Cell(13, 1).Range.Text = "18×55=990"
$t.Cell(13, 2).Range.Text = "100×76=7600"
$t.Cell(13, 3).Range.Text = "40×100=4000"
$t.Cell(13, 4).Range.Text = "96×40=3840"
$t.Cell(13, 5).Range.Text = "95×46=4370"
$t.Cell(14, 1).Range.Text = "74×86=6364"
$t.Cell(14, 2).Range.Text = "10×24=240"
$t.Cell(14, 3).Range.Text = "90×72=6480"
$t.Cell(14, 4).Range.Text = "28×64=1792"
$t.Cell(14, 5).Range.Text = "88×63=5544"
$t.Cell(15, 1).Range.Text = "58×49=2842"
$t.Cell(15, 2).Range.Text = "54×69=3726"
$t.Cell(15, 3).Range.Text = "79×56=4424"
$t.Cell(15, 4).Range.Text = "91×26=2366"
$t.Cell(15, 5).Range.Text = "54×72=3888"
$t.Cell(16, 1).Range.Text = "59×21=1239"
$t.Cell(16, 2).Range.Text = "17×14=238"
$t.Cell(16, 3).Range.Text = "54×43=2322"
$t.Cell(16, 4).Range.Text = "83×19=1577"
$t.Cell(16, 5).Range.Text = "82×42=3444"
$t.Cell(17, 1).Range.Text = "72×31=2232"
$t.Cell(17, 2).Range.Text = "54×64=3456"
$t.Cell(17, 3).Range.Text = "28×87=2436"
$t.Cell(17, 4).Range.Text = "74×37=2738"
$t.Cell(17, 5).Range.Text = "25×66=1650"
$t.Cell(18, 1).Range.Text = "30×92=2760"
$t.Cell(18, 2).Range.Text = "78×76=5928"
$t.Cell(18, 3).Range.Text = "77×88=6776"
$t.Cell(18, 4).Range.Text = "72×75=5400"
$t.Cell(18, 5).Range.Text = "96×52=4992"
$t.Cell(19, 1).Range.Text = "69×27=1863"
$t.Cell(19, 2).Range.Text = "72×51=3672"
$t.Cell(19, 3).Range.Text = "82×40=3280"
$t.Cell(19, 4).Range.Text = "56×21=1176"
$t.Cell(19, 5).Range.Text = "74×23=1702"
$t.Cell(20, 1).Range.Text = "40×84=3360"
$t.Cell(20, 2).Range.Text = "61×20=1220"
$t.Cell(20, 3).Range.Text = "42×19=798"
$t.Cell(20, 4).Range.Text = "85×69=5865"
$t.Cell(20, 5).Range.Text = "87×28=2436"
